$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Queries sheet: fix the typo in the SQL query string
# (" =, Last Changed On" -> "Last Changed On") and adjust the
# height of the row holding the query text.
# ---------------------------------------------------------------
$wsQueries = $wb.Worksheets.Item("Queries")
$query = @"
SELECT [TemplateName] as 'Template Name'
      ,[FileName] as 'FileName'
      ,[LastChangedBy] as 'Last Changed By'
  ,Substring(LastChangedOn,11,18) as 'Last Changed On'
  FROM [Product_OCM].[dbo].[FaxTemplate];
"@
$wsQueries.Range("A2").Value2 = $query
$wsQueries.Rows.Item(2).RowHeight = 75

# ---------------------------------------------------------------
# Create sheet: replace the sample row data with the new values.
# ---------------------------------------------------------------
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("A2").Value = "'Mode"
$wsCreate.Range("C2").Value = "'Colors Group"
$wsCreate.Range("D2").Value = "'Mode Custom Template.html"
$wsCreate.Range("D2").Select()

# ---------------------------------------------------------------
# Edit sheet: same new sample row data, plus a page setup change.
# ---------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "'Mode"
$wsEdit.Range("B2").Value = "'2"
$wsEdit.Range("C2").Value = "'Colors Group"
$wsEdit.Range("D2").Value = "'Mode Custom Template.html"
$wsEdit.Range("E2").Value = "'Sachin Score"
$wsEdit.PageSetup.PaperSize = 9
$wsEdit.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# Delete sheet: same new sample row data; this sheet becomes the
# selected / active tab.
# ---------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("A2").Value = "'Mode"
$wsDelete.Range("B2").Value = "'2"
$wsDelete.Range("C2").Value = "'Colors Group"
$wsDelete.Range("D2").Value = "'Mode Custom Template.html"
$wsDelete.Range("E2").Value = "'Sachin Score"
$wsDelete.Select()
